$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 9: new Antibiotic_name / Antibiotic_FinalConcentration for a higher
#     (200 ug/mL) concentration entry. ---
$ws.Range("O9").Value = "Ampicillin Antibiotic"
$ws.Range("P9").Value = "200_ug/mL"

# --- Row 7: new Antibiotic_name / Antibiotic_FinalConcentration for a
#     differently-formatted (no underscore) concentration entry. ---
$ws.Range("O7").Value = "Ampicillin Antibiotic"
$ws.Range("P7").Value = "0.1ug/mL"

# --- Row 9: corrected reagents/temperature/duration JSON in the Options
#     column (the original had a missing closing brace and an unquoted
#     "C" unit value). ---
$ws.Range("Q9").Value = '{"reagents": {"Ethanol": {"qty": 70, "units": "percent"} }, "temperature": {"qty": 37, "units": "C"}, "duration": {"qty": 15, "units": "minute"}}'

# --- Row 2: Antibiotic_name / Antibiotic_FinalConcentration are unchanged
#     text, but the Options column gets the same corrected JSON as above. ---
$ws.Range("O2").Value = "Ampicillin Antibiotic"
$ws.Range("P2").Value = "0.1_ug/mL"
$ws.Range("Q2").Value = '{"reagents": {"Ethanol": {"qty": 70, "units": "percent"} }, "temperature": {"qty": 37, "units": "C"}, "duration": {"qty": 15, "units": "minute"}}'

# --- View state: scroll right so column J is the left-most visible column
#     and leave the selection on Q9 (the last cell touched). ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("Q9").Select()
